$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "25.755.90"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.750.44"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.44%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "235.96"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.85%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.15%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5059"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.31%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 2)
$c.NumberFormat = "@"
$c.Value = "Cardano"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2716"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  +13.60%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 2)
$c.NumberFormat = "@"
$c.Value = "Dogecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.06211"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.28%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedEther"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "1.751.09"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 2)
$c.NumberFormat = "@"
$c.Value = "TRON"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.06925"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.17%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 2)
$c.NumberFormat = "@"
$c.Value = "Solana"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "15.53"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.24%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 2)
$c.NumberFormat = "@"
$c.Value = "Polygon"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.6127"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.68%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "78.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.46%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = "@"
$c.Value = "Polkadot"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "4.481"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.57%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 2)
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 2)
$c.NumberFormat = "@"
$c.Value = "Dai"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.10%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "25.778.48"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.20%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 2)
$c.NumberFormat = "@"
$c.Value = "Avalanche"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "11.66"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.82%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 2)
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "0.000006723"
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.36%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "1.978.47"
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.93%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 2)
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "4.045"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.19%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 2)
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "8.226"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.64%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 2)
$c.NumberFormat = "@"
$c.Value = "Chainlink"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "5.173"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.54%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "137.05"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.34%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = "@"
$c.Value = "Toncoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "1.456"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.06%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "15.13"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.19%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 2)
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "1.787"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.30%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 2)
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "102.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.64%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 2)
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.08279"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.98%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 2)
$c.NumberFormat = "@"
$c.Value = "InternetComputer(DFINITY)"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.729"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 2)
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "3.425"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.38%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 2)
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.04375"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.35%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 2)
$c.NumberFormat = "@"
$c.Value = "Frax"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.9991"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 2)
$c.NumberFormat = "@"
$c.Value = "HuobiToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "2.648"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.27%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 2)
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.91%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.6032"
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.59%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.679"
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.56%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "1.950"
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  -6.25%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 2)
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.01552"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.63%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = "@"
$c.Value = "PaxDollar"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 2)
$c.NumberFormat = "@"
$c.Value = "Quant"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "102.03"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.98%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 2)
$c.NumberFormat = "@"
$c.Value = "TrustWalletToken"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.7507"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.88%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 2)
$c.NumberFormat = "@"
$c.Value = "TheSandbox"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.3808"
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.14%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 2)
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "4.842"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.66%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = "Cronos"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.05497"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.60%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 2)
$c.NumberFormat = "@"
$c.Value = "Algorand"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.1084"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.08%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = "Elrond"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "30.28"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.67%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = "Aptos"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "5.940"
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.72%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = "Aave"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "52.28"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.70%  "
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = "@"
$c.Value = "USDD"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.24%  "
$c.Style = "Normal"
